# Scheduled market-data refresh: overwrite the cached FFXIV Marketboard
# price/profit figures (columns H-N) on each job sheet with freshly
# pulled Universalis averages. Mirrors the nightly runner diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 957.1177
$ws.Range("I43").Value = 953.5
$ws.Range("J43").Value = 959.0909
$ws.Range("K43").Value = 953.5
$ws.Range("L43").Value = 959.0909
$ws.Range("M43").Value = -884.5
$ws.Range("N43").Value = -1097.0909
# Row 76
$ws.Range("H76").Value = 9563.55
$ws.Range("I76").Value = 15877.1
$ws.Range("J76").Value = 3250
$ws.Range("K76").Value = 15877.1
$ws.Range("L76").Value = 3250
$ws.Range("M76").Value = -15562.1
$ws.Range("N76").Value = -3880
# Row 79
$ws.Range("H79").Value = 9563.55
$ws.Range("I79").Value = 15877.1
$ws.Range("J79").Value = 3250
$ws.Range("K79").Value = 15877.1
$ws.Range("L79").Value = 3250
$ws.Range("M79").Value = -14785.1
$ws.Range("N79").Value = -5434
# Row 116
$ws.Range("H116").Value = 4515.5
$ws.Range("I116").Value = 3388.3333
$ws.Range("J116").Value = 5360.875
$ws.Range("K116").Value = 3388.3333
$ws.Range("L116").Value = 5360.875
$ws.Range("M116").Value = 53.66670000000022
$ws.Range("N116").Value = -12244.875

$ws = $wb.Worksheets.Item("ARM")
# Row 44
$ws.Range("H44").Value = 34849
$ws.Range("J44").Value = 34849
$ws.Range("L44").Value = 34849
$ws.Range("N44").Value = -35825
# Row 55
$ws.Range("H55").Value = 33853
$ws.Range("J55").Value = 33853
$ws.Range("L55").Value = 33853
$ws.Range("N55").Value = -34483
# Row 63
$ws.Range("H63").Value = 3432.3684
$ws.Range("I63").Value = 2313.125
$ws.Range("J63").Value = 4246.364
$ws.Range("K63").Value = 2313.125
$ws.Range("L63").Value = 4246.364
$ws.Range("M63").Value = -1627.125
$ws.Range("N63").Value = -5618.364
# Row 66
$ws.Range("H66").Value = 3432.3684
$ws.Range("I66").Value = 2313.125
$ws.Range("J66").Value = 4246.364
$ws.Range("K66").Value = 11565.625
$ws.Range("L66").Value = 21231.82
$ws.Range("M66").Value = -8133.625
$ws.Range("N66").Value = -28095.82
# Row 80
$ws.Range("H80").Value = 42349.5
$ws.Range("J80").Value = 42349.5
$ws.Range("L80").Value = 42349.5
$ws.Range("N80").Value = -44345.5
# Row 83
$ws.Range("H83").Value = 42349.5
$ws.Range("J83").Value = 42349.5
$ws.Range("L83").Value = 127048.5
$ws.Range("N83").Value = -137032.5

$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Range("H35").Value = 34987
$ws.Range("J35").Value = 34987
$ws.Range("L35").Value = 34987
$ws.Range("N35").Value = -35607
# Row 82
$ws.Range("H82").Value = 32167.875
$ws.Range("J82").Value = 34634.715
$ws.Range("L82").Value = 34634.715
$ws.Range("N82").Value = -35400.715
# Row 85
$ws.Range("H85").Value = 32167.875
$ws.Range("J85").Value = 34634.715
$ws.Range("L85").Value = 34634.715
$ws.Range("N85").Value = -37286.715

$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 16914
$ws.Range("J41").Value = 19917.5
$ws.Range("L41").Value = 19917.5
$ws.Range("N41").Value = -20773.5
# Row 55
$ws.Range("H55").Value = 12321.5
$ws.Range("J55").Value = 13386
$ws.Range("L55").Value = 13386
$ws.Range("N55").Value = -14016
# Row 62
$ws.Range("H62").Value = 3931.25
$ws.Range("J62").Value = 4000
$ws.Range("L62").Value = 4000
$ws.Range("N62").Value = -5248
# Row 65
$ws.Range("H65").Value = 3931.25
$ws.Range("J65").Value = 4000
$ws.Range("L65").Value = 20000
$ws.Range("N65").Value = -26240
# Row 68
$ws.Range("H68").Value = 16903.6
$ws.Range("J68").Value = 16903.6
$ws.Range("L68").Value = 16903.6
$ws.Range("N68").Value = -18401.6
# Row 71
$ws.Range("H71").Value = 16903.6
$ws.Range("J71").Value = 16903.6
$ws.Range("L71").Value = 50710.8
$ws.Range("N71").Value = -58198.8
# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 357.6154
$ws.Range("I5").Value = 339.13043
$ws.Range("J5").Value = 499.33334
$ws.Range("K5").Value = 1017.39129
$ws.Range("L5").Value = 1498.00002
$ws.Range("M5").Value = -905.39129
$ws.Range("N5").Value = -1722.00002
# Row 131
$ws.Range("H131").Value = 885.2347
$ws.Range("I131").Value = 565.7143
$ws.Range("J131").Value = 938.4881
$ws.Range("K131").Value = 1697.1429
$ws.Range("L131").Value = 2815.4643
$ws.Range("M131").Value = 3342.8571
$ws.Range("N131").Value = -12895.4643
# Row 135
$ws.Range("H135").Value = 357.6154
$ws.Range("I135").Value = 339.13043
$ws.Range("J135").Value = 499.33334
$ws.Range("K135").Value = 3052.17387
$ws.Range("L135").Value = 4494.00006
$ws.Range("M135").Value = -517.1738700000001
$ws.Range("N135").Value = -9564.00006
# Row 136
$ws.Range("H136").Value = 3341.9673
$ws.Range("I136").Value = 1341.6666
$ws.Range("J136").Value = 4179.3022
$ws.Range("K136").Value = 4024.9998
$ws.Range("L136").Value = 12537.9066
$ws.Range("M136").Value = 1075.0002
$ws.Range("N136").Value = -22737.9066

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 5157
$ws.Range("J43").Value = 8349.75
$ws.Range("L43").Value = 8349.75
$ws.Range("N43").Value = -8651.75
# Row 80
$ws.Range("H80").Value = 3135.625
$ws.Range("I80").Value = 2817
$ws.Range("J80").Value = 3666.6667
$ws.Range("K80").Value = 2817
$ws.Range("L80").Value = 3666.6667
$ws.Range("M80").Value = -1819
$ws.Range("N80").Value = -5662.6667
# Row 83
$ws.Range("H83").Value = 3135.625
$ws.Range("I83").Value = 2817
$ws.Range("J83").Value = 3666.6667
$ws.Range("K83").Value = 14085
$ws.Range("L83").Value = 18333.3335
$ws.Range("M83").Value = -9093
$ws.Range("N83").Value = -28317.3335
# Row 102
$ws.Range("H102").Value = 2878
$ws.Range("I102").Value = 1375.5
$ws.Range("J102").Value = 4080
$ws.Range("K102").Value = 1375.5
$ws.Range("L102").Value = 4080
$ws.Range("M102").Value = 246.5
$ws.Range("N102").Value = -7324
# Row 122
$ws.Range("H122").Value = 3001.4
$ws.Range("I122").Value = 3002.3333
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 9006.999899999999
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -6556.999899999999
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1190.6666
$ws.Range("I22").Value = 842.8571
$ws.Range("J22").Value = 1412
$ws.Range("K22").Value = 842.8571
$ws.Range("L22").Value = 1412
$ws.Range("M22").Value = -547.8571
$ws.Range("N22").Value = -2002
# Row 27
$ws.Range("H27").Value = 1190.6666
$ws.Range("I27").Value = 842.8571
$ws.Range("J27").Value = 1412
$ws.Range("K27").Value = 842.8571
$ws.Range("L27").Value = 1412
$ws.Range("M27").Value = -735.8571
$ws.Range("N27").Value = -1626
# Row 40
$ws.Range("H40").Value = 2070.3
$ws.Range("I40").Value = 2070.3
$ws.Range("K40").Value = 2070.3
$ws.Range("M40").Value = -1934.3
# Row 93
$ws.Range("H93").Value = 2149.6924
$ws.Range("I93").Value = 1894.6
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 1894.6
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -646.5999999999999
$ws.Range("N93").Value = -5496
# Row 109
$ws.Range("H109").Value = 22761.666
$ws.Range("J109").Value = 22761.666
$ws.Range("L109").Value = 22761.666
$ws.Range("N109").Value = -25535.666
# Row 122
$ws.Range("H122").Value = 3016.8262
$ws.Range("I122").Value = 2483.2307
$ws.Range("J122").Value = 3710.5
$ws.Range("K122").Value = 7449.6921
$ws.Range("L122").Value = 11131.5
$ws.Range("M122").Value = -4999.6921
$ws.Range("N122").Value = -16031.5

$ws = $wb.Worksheets.Item("WVR")
# Row 109
$ws.Range("H109").Value = 18788.5
$ws.Range("J109").Value = 18788.5
$ws.Range("L109").Value = 18788.5

# Rows whose market data dropped to zero lose their derived profit cell
# entirely (matches the sheet convention: no price => no profit figure).
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N109").ClearContents()
